# Chamada 3DES de 08-01-2021
# Adds the missing Friday (08-01-2021) attendance column (F) for every
# student row, plus a whole new week (columns G:K, 11-01-2021 to
# 15-01-2021) with its own header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New week header row (row 1): PROJ, PROJ, PROJ, PDMO, RMST ---
$ws.Range("G1").Value = "PROJ"
$ws.Range("H1").Value = "PROJ"
$ws.Range("I1").Value = "PROJ"
$ws.Range("J1").Value = "PDMO"
$ws.Range("K1").Value = "RMST"

# --- New week dates (row 2): 11-01-2021 .. 15-01-2021, same date format as F2 ---
$ws.Range("G2:K2").NumberFormat = $ws.Range("F2").NumberFormat
$ws.Range("G2").Value = 44207
$ws.Range("H2").Value = 44208
$ws.Range("I2").Value = 44209
$ws.Range("J2").Value = 44210
$ws.Range("K2").Value = 44211

# --- Fill in the previously-missing Friday (08-01-2021) attendance, column F ---
$ws.Range("F3").Value = "P"
$ws.Range("F4").Value = "P"
$ws.Range("F5").Value = "P"
$ws.Range("F6").Value = "P"
$ws.Range("F7").Value = "F"
$ws.Range("F8").Value = "F"
$ws.Range("F9").Value = "F"
$ws.Range("F10").Value = "F"
$ws.Range("F11").Value = "P"
$ws.Range("F12").Value = "P"
$ws.Range("F13").Value = "P"
$ws.Range("F14").Value = "P"
$ws.Range("F15").Value = "P"
$ws.Range("F16").Value = "P"
$ws.Range("F17").Value = "P"
$ws.Range("F18").Value = "P"
$ws.Range("F19").Value = "P"
$ws.Range("F20").Value = "P"

# --- Selection ends up on J4, matching the author's last click ---
[void]$ws.Range("J4").Select()
